$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The edit merges the first two paragraphs ("Una nube sobre el agua "
# + "Cae la lluvia y el pasto") into a single paragraph and replaces
# all of its text with a new sentence (with a couple of spell-check
# markers around "we" and "se", as Word would leave behind), keeping
# the trailing "_GoBack" bookmark at the very end of the paragraph.
# ------------------------------------------------------------------

# 0) Drop the existing "_GoBack" bookmark up front - we recreate it
#    (collapsed, right after the new text) as part of the XML we
#    insert below, so it doesn't get stretched across the whole
#    paragraph by the edits that follow.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 1) Merge paragraph 1 and paragraph 2 by deleting the paragraph mark
#    that ends paragraph 1 (equivalent to placing the cursor at the
#    end of paragraph 1 and pressing Delete in Word).
$p1 = $d.Paragraphs(1)
$paraMark = $d.Range($p1.Range.End - 1, $p1.Range.End)
$paraMark.Delete()

# 2) Replace the merged paragraph's text (everything up to, but not
#    including, its own trailing paragraph mark) with the new runs.
$p = $d.Paragraphs(1)
$body = $d.Range($p.Range.Start, $p.Range.End - 1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:lang w:val="es-CL"/></w:rPr>'

$xml = "<w:p $wNs>" + `
    "<w:r>$rPr<w:t xml:space=`"preserve`">Esta </w:t></w:r>" + `
    '<w:proofErr w:type="spellStart"/>' + `
    "<w:r>$rPr<w:t>we</w:t></w:r>" + `
    '<w:proofErr w:type="spellEnd"/>' + `
    "<w:r>$rPr<w:t xml:space=`"preserve`"> no funciona y no </w:t></w:r>" + `
    '<w:proofErr w:type="spellStart"/>' + `
    "<w:r>$rPr<w:t>se</w:t></w:r>" + `
    '<w:proofErr w:type="spellEnd"/>' + `
    "<w:r>$rPr<w:t xml:space=`"preserve`"> porque ¡!!!</w:t></w:r>" + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'

$body.InsertXML($xml)
